$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Merge the two bullet paragraphs into one: replace the paragraph mark
# between the two sentences with a single space so they become one
# paragraph / one run of text.
$searchText  = "Phonegap." + [char]13 + "Collaborated"
$replaceText = "Phonegap. Collaborated"

$find.Execute($searchText, $true, $false, $false, $false, $false, `
              $true, 1, $false, $replaceText, 2) | Out-Null
